$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin (B), Link (C), Price (D), Volume1h (E)
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.365.29", "  +1.54%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.010.71", "  +4.84%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.002", "  -0.02%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "324.72", "  +1.35%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  -0.04%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5129", "  +1.41%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.4263", "  +5.98%  "),
    @(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.08761", "  +5.38%  "),
    @(10, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.134", "  +2.92%  "),
    @(11, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "24.48", "  +3.14%  "),
    @(12, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.013.80", "  +5.15%  "),
    @(13, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "6.607", "  +3.31%  "),
    @(14, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.458", "  +3.23%  "),
    @(15, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.003", "  +0.31%  "),
    @(16, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "94.26", "  +2.27%  "),
    @(17, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001114", "  +1.51%  "),
    @(18, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06538", "  +0.36%  "),
    @(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "18.85", "  +3.62%  "),
    @(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  +0.04%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.205", "  +4.59%  "),
    @(22, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.417.77", "  +1.64%  "),
    @(23, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.85", "  +4.80%  "),
    @(24, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.263", "  +3.10%  "),
    @(25, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.253.05", "  +5.59%  "),
    @(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "22.45", "  +1.77%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "162.53", "  +0.22%  "),
    @(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.435", "  +5.07%  "),
    @(29, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "131.26", "  +1.78%  "),
    @(30, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.149", "  +1.90%  "),
    @(31, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1054", "  +1.73%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "6.105", "  +2.45%  "),
    @(33, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.833", "  +0.64%  "),
    @(34, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.363", "  +14.44%  "),
    @(35, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02529", "  +3.36%  "),
    @(36, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.473", "  +1.37%  "),
    @(37, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06670", "  +3.98%  "),
    @(38, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "12.43", "  +9.47%  "),
    @(39, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "9.171", "  +5.38%  "),
    @(40, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2215", "  +2.95%  "),
    @(41, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.6669", "  +2.69%  "),
    @(42, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.233", "  +1.61%  "),
    @(43, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "1.001", "  +0.03%  "),
    @(44, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "13.68", "  +2.44%  "),
    @(45, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.6184", "  +2.27%  "),
    @(46, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "2.197", "  -1.36%  "),
    @(47, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.624", "  -0.36%  "),
    @(48, "EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "1.258", "  +4.16%  "),
    @(49, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "124.86", "  +2.38%  "),
    @(50, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "81.22", "  +3.08%  "),
    @(51, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06915", "  +1.66%  ")

)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
}
